# Trade Request Example.xlsx - refactor of trade_request pricing, portfolio
#
# The "model" worksheet's portfolio/trade-request rows are refactored:
#   1. The stale "Almost Cash" placeholder row (old row 17: A="45-33",
#      B="Almost Cash", D=10000) is removed entirely, shifting every row
#      below it up by one.
#   2. The old "LQD" line in the 5268-5955 portfolio (old row 25, now row
#      24 after the shift) is repriced/retargeted: its symbol becomes
#      "SPY" and its stale price (E) is cleared, while its share count
#      (D=121.1987) is kept as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# 1. Remove the obsolete "Almost Cash" row - everything below shifts up.
$ws.Rows.Item(17).Delete()

# 2. Re-point the old "LQD" row (now row 24) at "SPY" and drop its price.
$ws.Range("B24").Value = "SPY"
$ws.Range("E24").ClearContents()

# Update the sheet view: drop the frozen/scrolled topLeftCell and move the
# active selection to E24.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E24").Select()
